$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Values that are valid numeric literals (e.g. "1.004")
# are prefixed with a leading apostrophe so Excel stores them as quote-prefixed
# TEXT (matching the workbook's original inlineStr/text cells) instead of
# silently converting them to numbers.
$updates = @{
    'D2' = '26.528.27'
    'E2' = '  -7.58%  '
    'D3' = '1.681.04'
    'E3' = '  -6.78%  '
    'D4' = '''1.004'
    'E4' = '  +0.16%  '
    'D5' = '''216.07'
    'E5' = '  -6.69%  '
    'D6' = '''1.004'
    'E6' = '  +0.10%  '
    'D7' = '''0.4980'
    'E7' = '  -16.25%  '
    'D8' = '''0.2601'
    'E8' = '  -6.57%  '
    'D9' = '''21.82'
    'E9' = '  -6.54%  '
    'D10' = '''0.06191'
    'E10' = '  -9.49%  '
    'D11' = '''0.07280'
    'E11' = '  -3.50%  '
    'D12' = '1.676.54'
    'E12' = '  -4.66%  '
    'D13' = '''4.433'
    'E13' = '  -7.21%  '
    'D14' = '''0.5744'
    'E14' = '  -7.93%  '
    'D15' = '1.909.08'
    'E15' = '  -6.81%  '
    'D16' = '''0.000008169'
    'E16' = '  -12.68%  '
    'D17' = '''64.39'
    'E17' = '  -14.81%  '
    'D18' = '26.522.97'
    'E18' = '  -7.46%  '
    'D19' = '''4.986'
    'E19' = '  -9.21%  '
    'E20' = '  +0.15%  '
    'D21' = '''10.76'
    'E21' = '  -6.07%  '
    'D22' = '''184.35'
    'E22' = '  -12.31%  '
    'D23' = '''6.182'
    'E23' = '  -9.89%  '
    'D24' = '''1.005'
    'E24' = '  +0.19%  '
    'D25' = '''144.40'
    'E25' = '  -6.37%  '
    'D26' = '''7.442'
    'E26' = '  -5.42%  '
    'D27' = '''0.1130'
    'E27' = '  -11.19%  '
    'D28' = '''15.41'
    'E28' = '  -5.99%  '
    'D29' = '''1.302'
    'E29' = '  -8.90%  '
    'D30' = '''0.05693'
    'E30' = '  -8.08%  '
    'D31' = '''1.318'
    'E31' = '  -7.22%  '
    'D32' = '''3.473'
    'E32' = '  -8.18%  '
    'D33' = '''3.463'
    'E33' = '  -7.66%  '
    'D34' = '''1.632'
    'E34' = '  -5.16%  '
    'D35' = '''1.005'
    'E35' = '  -5.68%  '
    'E36' = '  -5.01%  '
    'D37' = '''0.5905'
    'E37' = '  -7.75%  '
    'E38' = '  -3.02%  '
    'D39' = '''0.01588'
    'E39' = '  -7.50%  '
    'D40' = '1.068.49'
    'E40' = '  -5.69%  '
    'D41' = '''5.870'
    'E41' = '  -9.21%  '
    'D42' = '''0.8518'
    'E42' = '  -2.74%  '
    'D43' = '''1.002'
    'E43' = '  -0.39%  '
    'D44' = '''98.22'
    'E44' = '  -2.50%  '
    'D45' = '1.835.88'
    'E45' = '  -6.41%  '
    'D46' = '''56.13'
    'E46' = '  -7.30%  '
    'B47' = 'BabyDogeCoin'
    'C47' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D47' = '''0.00000000105'
    'E47' = '  -5.96%  '
    'B48' = 'Frax'
    'C48' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'D48' = '''1.003'
    'E48' = '  -0.19%  '
    'D49' = '''8.063'
    'E49' = '  -3.42%  '
    'D50' = '''0.4311'
    'E50' = '  -3.88%  '
    'D51' = '''0.05188'
    'E51' = '  -5.23%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
